$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update column F (想去人数 / interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 782
$ws1.Range("F6").Value = 67
$ws1.Range("F7").Value = 271
$ws1.Range("F8").Value = 3880
$ws1.Range("F9").Value = 87
$ws1.Range("F10").Value = 4572
$ws1.Range("F12").Value = 1151
$ws1.Range("F13").Value = 71

# Sheet "全部类型" (all types) - update column F (想去人数 / interested count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 782
$ws4.Range("F6").Value = 67
$ws4.Range("F8").Value = 271
$ws4.Range("F9").Value = 3880
$ws4.Range("F10").Value = 87
$ws4.Range("F11").Value = 4572
$ws4.Range("F13").Value = 1151
$ws4.Range("F14").Value = 71
